$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$lastRow = 56
for ($r = 2; $r -le $lastRow; $r++) {
    $cell = $ws.Cells.Item($r, 3)
    if ($cell.Value2 -eq 46074) {
        $cell.Value = 46075
    }
}
